$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 461  # ALC H12 (was 458.85715)
$ws.Cells.Item(12, 9).Value = 475  # ALC I12 (was 474.75)
$ws.Cells.Item(12, 11).Value = 475  # ALC K12 (was 474.75)
$ws.Cells.Item(12, 13).Value = -305  # ALC M12 (was -304.75)

$ws.Cells.Item(17, 8).Value = 961.3393  # ALC H17 (was 991.5848999999999)
$ws.Cells.Item(17, 10).Value = 969.9636  # ALC J17 (was 1001.28845)
$ws.Cells.Item(17, 12).Value = 2909.8908  # ALC L17 (was 3003.86535)
$ws.Cells.Item(17, 14).Value = -3245.8908  # ALC N17 (was -3339.86535)

$ws.Cells.Item(48, 8).Value = 1978.0625  # ALC H48 (was 1953)
$ws.Cells.Item(48, 9).Value = 1939.9286  # ALC I48 (was 2012.6154)
$ws.Cells.Item(48, 10).Value = 2245  # ALC J48 (was 1798)
$ws.Cells.Item(48, 11).Value = 5819.7858  # ALC K48 (was 6037.8462)
$ws.Cells.Item(48, 12).Value = 6735  # ALC L48 (was 5394)
$ws.Cells.Item(48, 13).Value = -5527.7858  # ALC M48 (was -5745.8462)
$ws.Cells.Item(48, 14).Value = -7319  # ALC N48 (was -5978)

$ws.Cells.Item(55, 8).Value = 612.36365  # ALC H55 (was 672.6)
$ws.Cells.Item(55, 9).Value = 59.25  # ALC I55 (was 75.666664)
$ws.Cells.Item(55, 11).Value = 59.25  # ALC K55 (was 75.666664)
$ws.Cells.Item(55, 13).Value = 154.75  # ALC M55 (was 138.333336)

$ws.Cells.Item(56, 8).Value = 1978.0625  # ALC H56 (was 1953)
$ws.Cells.Item(56, 9).Value = 1939.9286  # ALC I56 (was 2012.6154)
$ws.Cells.Item(56, 10).Value = 2245  # ALC J56 (was 1798)
$ws.Cells.Item(56, 11).Value = 5819.7858  # ALC K56 (was 6037.8462)
$ws.Cells.Item(56, 12).Value = 6735  # ALC L56 (was 5394)
$ws.Cells.Item(56, 13).Value = -5285.7858  # ALC M56 (was -5503.8462)
$ws.Cells.Item(56, 14).Value = -7803  # ALC N56 (was -6462)

$ws.Cells.Item(70, 8).Value = 2525.923  # ALC H70 (was 2699.2727)
$ws.Cells.Item(70, 9).Value = 1597  # ALC I70 (was 1599)
$ws.Cells.Item(70, 10).Value = 2694.818  # ALC J70 (was 2809.3)
$ws.Cells.Item(70, 11).Value = 4791  # ALC K70 (was 4797)
$ws.Cells.Item(70, 12).Value = 8084.454000000001  # ALC L70 (was 8427.900000000001)
$ws.Cells.Item(70, 13).Value = -4521  # ALC M70 (was -4527)
$ws.Cells.Item(70, 14).Value = -8624.454000000002  # ALC N70 (was -8967.900000000001)

$ws.Cells.Item(73, 8).Value = 2525.923  # ALC H73 (was 2699.2727)
$ws.Cells.Item(73, 9).Value = 1597  # ALC I73 (was 1599)
$ws.Cells.Item(73, 10).Value = 2694.818  # ALC J73 (was 2809.3)
$ws.Cells.Item(73, 11).Value = 4791  # ALC K73 (was 4797)
$ws.Cells.Item(73, 12).Value = 8084.454000000001  # ALC L73 (was 8427.900000000001)
$ws.Cells.Item(73, 13).Value = -3855  # ALC M73 (was -3861)
$ws.Cells.Item(73, 14).Value = -9956.454000000002  # ALC N73 (was -10299.9)

$ws.Cells.Item(74, 8).Value = 14846.3  # ALC H74 (was 15885.889)
$ws.Cells.Item(74, 10).Value = 24393.6  # ALC J74 (was 29119.5)
$ws.Cells.Item(74, 12).Value = 24393.6  # ALC L74 (was 29119.5)
$ws.Cells.Item(74, 14).Value = -26265.6  # ALC N74 (was -30991.5)

$ws.Cells.Item(76, 8).Value = 5853.7144  # ALC H76 (was 5372)
$ws.Cells.Item(76, 9).Value = 4796.2  # ALC I76 (was 4330.1665)
$ws.Cells.Item(76, 11).Value = 4796.2  # ALC K76 (was 4330.1665)
$ws.Cells.Item(76, 13).Value = -4481.2  # ALC M76 (was -4015.1665)

$ws.Cells.Item(77, 8).Value = 14846.3  # ALC H77 (was 15885.889)
$ws.Cells.Item(77, 10).Value = 24393.6  # ALC J77 (was 29119.5)
$ws.Cells.Item(77, 12).Value = 121968  # ALC L77 (was 145597.5)
$ws.Cells.Item(77, 14).Value = -131328  # ALC N77 (was -154957.5)

$ws.Cells.Item(79, 8).Value = 5853.7144  # ALC H79 (was 5372)
$ws.Cells.Item(79, 9).Value = 4796.2  # ALC I79 (was 4330.1665)
$ws.Cells.Item(79, 11).Value = 4796.2  # ALC K79 (was 4330.1665)
$ws.Cells.Item(79, 13).Value = -3704.2  # ALC M79 (was -3238.1665)

$ws.Cells.Item(94, 8).Value = 5262.6665  # ALC H94 (was 5010.5713)
$ws.Cells.Item(94, 10).Value = 7000  # ALC J94 (was 5249)
$ws.Cells.Item(94, 12).Value = 7000  # ALC L94 (was 5249)
$ws.Cells.Item(94, 14).Value = -7902  # ALC N94 (was -6151)

$ws.Cells.Item(98, 8).Value = 2560.2258  # ALC H98 (was 2491.9614)
$ws.Cells.Item(98, 9).Value = 2512.2666  # ALC I98 (was 2431.68)
$ws.Cells.Item(98, 11).Value = 2512.2666  # ALC K98 (was 2431.68)
$ws.Cells.Item(98, 13).Value = -1014.2666  # ALC M98 (was -933.6799999999998)

$ws.Cells.Item(122, 8).Value = 2560.2258  # ALC H122 (was 2491.9614)
$ws.Cells.Item(122, 9).Value = 2512.2666  # ALC I122 (was 2431.68)
$ws.Cells.Item(122, 11).Value = 7536.7998  # ALC K122 (was 7295.039999999999)
$ws.Cells.Item(122, 13).Value = -5086.7998  # ALC M122 (was -4845.039999999999)

$ws.Cells.Item(129, 8).Value = 1252.3889  # ALC H129 (was 1251.7916)

$ws.Cells.Item(131, 8).Value = 16665  # ALC H131 (was 14997.5)
$ws.Cells.Item(131, 10).Value = 17500  # ALC J131 (was 15000)
$ws.Cells.Item(131, 12).Value = 52500  # ALC L131 (was 45000)
$ws.Cells.Item(131, 14).Value = -62580  # ALC N131 (was -55080)

$ws.Cells.Item(137, 8).Value = 50876.332  # ALC H137 (was 48521.547)
$ws.Cells.Item(137, 9).Value = 79555.16  # ALC I137 (was 73913.42999999999)
$ws.Cells.Item(137, 10).Value = 4273.25  # ALC J137 (was 4085.75)
$ws.Cells.Item(137, 11).Value = 238665.48  # ALC K137 (was 221740.29)
$ws.Cells.Item(137, 12).Value = 12819.75  # ALC L137 (was 12257.25)
$ws.Cells.Item(137, 13).Value = -236115.48  # ALC M137 (was -219190.29)
$ws.Cells.Item(137, 14).Value = -17919.75  # ALC N137 (was -17357.25)

$ws.Cells.Item(138, 8).Value = 16852.68  # ALC H138 (was 18153.152)
$ws.Cells.Item(138, 9).Value = 20225.883  # ALC I138 (was 20826.908)
$ws.Cells.Item(138, 10).Value = 9684.625  # ALC J138 (was 11365.923)
$ws.Cells.Item(138, 11).Value = 60677.649  # ALC K138 (was 62480.724)
$ws.Cells.Item(138, 12).Value = 29053.875  # ALC L138 (was 34097.769)
$ws.Cells.Item(138, 13).Value = -55537.649  # ALC M138 (was -57340.724)
$ws.Cells.Item(138, 14).Value = -39333.875  # ALC N138 (was -44377.769)

$ws.Cells.Item(141, 8).Value = 1568.8572  # ALC H141 (was 1637.2)
$ws.Cells.Item(141, 9).Value = 1517  # ALC I141 (was 1621.75)
$ws.Cells.Item(141, 10).Value = 1698.5  # ALC J141 (was 1699)
$ws.Cells.Item(141, 11).Value = 4551  # ALC K141 (was 4865.25)
$ws.Cells.Item(141, 12).Value = 5095.5  # ALC L141 (was 5097)
$ws.Cells.Item(141, 13).Value = 629  # ALC M141 (was 314.75)
$ws.Cells.Item(141, 14).Value = -15455.5  # ALC N141 (was -15457)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 30560.375  # ARM H32 (was 34588.855)
$ws.Cells.Item(32, 9).Value = 33380.03  # ARM I32 (was 38728.137)
$ws.Cells.Item(32, 11).Value = 33380.03  # ARM K32 (was 38728.137)
$ws.Cells.Item(32, 13).Value = -33093.03  # ARM M32 (was -38441.137)

$ws.Cells.Item(43, 8).Value = 19844  # ARM H43 (was 21591.25)
$ws.Cells.Item(43, 10).Value = 19844  # ARM J43 (was 21591.25)
$ws.Cells.Item(43, 12).Value = 19844  # ARM L43 (was 21591.25)
$ws.Cells.Item(43, 14).Value = -20470  # ARM N43 (was -22217.25)

$ws.Cells.Item(45, 8).Value = 3527.5293  # ARM H45 (was 3392.611)
$ws.Cells.Item(45, 9).Value = 3366.3333  # ARM I45 (was 3399.6667)
$ws.Cells.Item(45, 10).Value = 3562.0715  # ARM J45 (was 3391.2)
$ws.Cells.Item(45, 11).Value = 3366.3333  # ARM K45 (was 3399.6667)
$ws.Cells.Item(45, 12).Value = 3562.0715  # ARM L45 (was 3391.2)
$ws.Cells.Item(45, 13).Value = -2989.3333  # ARM M45 (was -3022.6667)
$ws.Cells.Item(45, 14).Value = -4316.0715  # ARM N45 (was -4145.2)

$ws.Cells.Item(61, 8).Value = 6269.6  # ARM H61 (was 8242.429)
$ws.Cells.Item(61, 9).Value = 4671  # ARM I61 (was 6924.5)
$ws.Cells.Item(61, 11).Value = 4671  # ARM K61 (was 6924.5)
$ws.Cells.Item(61, 13).Value = -4459  # ARM M61 (was -6712.5)

$ws.Cells.Item(74, 8).Value = 1736.75  # ARM H74 (was 1770.5714)
$ws.Cells.Item(74, 10).Value = 1330.6666  # ARM J74 (was 1246)
$ws.Cells.Item(74, 12).Value = 1330.6666  # ARM L74 (was 1246)
$ws.Cells.Item(74, 14).Value = -3078.6666  # ARM N74 (was -2994)

$ws.Cells.Item(76, 8).Value = 0  # ARM H76 (was 7396.4)
$ws.Cells.Item(76, 9).Value = 0  # ARM I76 (was 2800)
$ws.Cells.Item(76, 10).Value = 0  # ARM J76 (was 8545.5)
$ws.Cells.Item(76, 11).Value = 0  # ARM K76 (was 2800)
$ws.Cells.Item(76, 12).Value = 0  # ARM L76 (was 8545.5)
$ws.Cells.Item(76, 13).Value = $null  # ARM M76 remove (was -2462)
$ws.Cells.Item(76, 14).Value = $null  # ARM N76 remove (was -9221.5)

$ws.Cells.Item(77, 8).Value = 1736.75  # ARM H77 (was 1770.5714)
$ws.Cells.Item(77, 10).Value = 1330.6666  # ARM J77 (was 1246)
$ws.Cells.Item(77, 12).Value = 6653.333000000001  # ARM L77 (was 6230)
$ws.Cells.Item(77, 14).Value = -15389.333  # ARM N77 (was -14966)

$ws.Cells.Item(79, 8).Value = 0  # ARM H79 (was 7396.4)
$ws.Cells.Item(79, 9).Value = 0  # ARM I79 (was 2800)
$ws.Cells.Item(79, 10).Value = 0  # ARM J79 (was 8545.5)
$ws.Cells.Item(79, 11).Value = 0  # ARM K79 (was 2800)
$ws.Cells.Item(79, 12).Value = 0  # ARM L79 (was 8545.5)
$ws.Cells.Item(79, 13).Value = $null  # ARM M79 remove (was -1630)
$ws.Cells.Item(79, 14).Value = $null  # ARM N79 remove (was -10885.5)

$ws.Cells.Item(88, 8).Value = 3469.9167  # ARM H88 (was 3602.6365)
$ws.Cells.Item(88, 9).Value = 2920.5  # ARM I88 (was 3006)
$ws.Cells.Item(88, 10).Value = 3579.8  # ARM J88 (was 3662.3)
$ws.Cells.Item(88, 11).Value = 2920.5  # ARM K88 (was 3006)
$ws.Cells.Item(88, 12).Value = 3579.8  # ARM L88 (was 3662.3)
$ws.Cells.Item(88, 13).Value = -2514.5  # ARM M88 (was -2600)
$ws.Cells.Item(88, 14).Value = -4391.8  # ARM N88 (was -4474.3)

$ws.Cells.Item(91, 8).Value = 3469.9167  # ARM H91 (was 3602.6365)
$ws.Cells.Item(91, 9).Value = 2920.5  # ARM I91 (was 3006)
$ws.Cells.Item(91, 10).Value = 3579.8  # ARM J91 (was 3662.3)
$ws.Cells.Item(91, 11).Value = 2920.5  # ARM K91 (was 3006)
$ws.Cells.Item(91, 12).Value = 3579.8  # ARM L91 (was 3662.3)
$ws.Cells.Item(91, 13).Value = -1516.5  # ARM M91 (was -1602)
$ws.Cells.Item(91, 14).Value = -6387.8  # ARM N91 (was -6470.3)

$ws.Cells.Item(97, 8).Value = 11832.363  # ARM H97 (was 15674)
$ws.Cells.Item(97, 9).Value = 17629.334  # ARM I97 (was 34265.668)
$ws.Cells.Item(97, 10).Value = 4876  # ARM J97 (was 4519)
$ws.Cells.Item(97, 11).Value = 17629.334  # ARM K97 (was 34265.668)
$ws.Cells.Item(97, 12).Value = 4876  # ARM L97 (was 4519)
$ws.Cells.Item(97, 13).Value = -17133.334  # ARM M97 (was -33769.668)
$ws.Cells.Item(97, 14).Value = -5868  # ARM N97 (was -5511)

$ws.Cells.Item(122, 8).Value = 1278.5186  # ARM H122 (was 1341.84)
$ws.Cells.Item(122, 9).Value = 1232.7693  # ARM I122 (was 1294.9166)
$ws.Cells.Item(122, 11).Value = 3698.3079  # ARM K122 (was 3884.7498)
$ws.Cells.Item(122, 13).Value = -1248.3079  # ARM M122 (was -1434.7498)

$ws.Cells.Item(132, 8).Value = 26392  # ARM H132 (was 28256.1)
$ws.Cells.Item(132, 9).Value = 28969.947  # ARM I132 (was 31321.314)
$ws.Cells.Item(132, 11).Value = 86909.841  # ARM K132 (was 93963.942)
$ws.Cells.Item(132, 13).Value = -84379.841  # ARM M132 (was -91433.942)

$ws.Cells.Item(135, 8).Value = 50000  # ARM H135 (was 55000)
$ws.Cells.Item(135, 10).Value = 50000  # ARM J135 (was 55000)
$ws.Cells.Item(135, 12).Value = 50000  # ARM L135 (was 55000)
$ws.Cells.Item(135, 14).Value = -60140  # ARM N135 (was -65140)

$ws.Cells.Item(136, 8).Value = 6269.6  # ARM H136 (was 8242.429)
$ws.Cells.Item(136, 9).Value = 4671  # ARM I136 (was 6924.5)
$ws.Cells.Item(136, 11).Value = 14013  # ARM K136 (was 20773.5)
$ws.Cells.Item(136, 13).Value = -11463  # ARM M136 (was -18223.5)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4303.4165  # BSM H86 (was 4094)
$ws.Cells.Item(86, 9).Value = 3411.25  # BSM I86 (was 3174.8572)
$ws.Cells.Item(86, 10).Value = 4749.5  # BSM J86 (was 5166.3335)
$ws.Cells.Item(86, 11).Value = 3411.25  # BSM K86 (was 3174.8572)
$ws.Cells.Item(86, 12).Value = 4749.5  # BSM L86 (was 5166.3335)
$ws.Cells.Item(86, 13).Value = -2288.25  # BSM M86 (was -2051.8572)
$ws.Cells.Item(86, 14).Value = -6995.5  # BSM N86 (was -7412.3335)

$ws.Cells.Item(89, 8).Value = 4303.4165  # BSM H89 (was 4094)
$ws.Cells.Item(89, 9).Value = 3411.25  # BSM I89 (was 3174.8572)
$ws.Cells.Item(89, 10).Value = 4749.5  # BSM J89 (was 5166.3335)
$ws.Cells.Item(89, 11).Value = 17056.25  # BSM K89 (was 15874.286)
$ws.Cells.Item(89, 12).Value = 23747.5  # BSM L89 (was 25831.6675)
$ws.Cells.Item(89, 13).Value = -11440.25  # BSM M89 (was -10258.286)
$ws.Cells.Item(89, 14).Value = -34979.5  # BSM N89 (was -37063.6675)

$ws.Cells.Item(94, 8).Value = 5015.1113  # BSM H94 (was 4140.077)
$ws.Cells.Item(94, 9).Value = 4424.75  # BSM I94 (was 3442)
$ws.Cells.Item(94, 10).Value = 5487.4  # BSM J94 (was 4954.5)
$ws.Cells.Item(94, 11).Value = 4424.75  # BSM K94 (was 3442)
$ws.Cells.Item(94, 12).Value = 5487.4  # BSM L94 (was 4954.5)
$ws.Cells.Item(94, 13).Value = -3973.75  # BSM M94 (was -2991)
$ws.Cells.Item(94, 14).Value = -6389.4  # BSM N94 (was -5856.5)

$ws.Cells.Item(105, 8).Value = 5340.5  # BSM H105 (was 5413.125)
$ws.Cells.Item(105, 9).Value = 5239.5  # BSM I105 (was 5334.25)
$ws.Cells.Item(105, 11).Value = 5239.5  # BSM K105 (was 5334.25)
$ws.Cells.Item(105, 13).Value = -3492.5  # BSM M105 (was -3587.25)

$ws.Cells.Item(134, 8).Value = 2226.9583  # BSM H134 (was 2366.5217)
$ws.Cells.Item(134, 9).Value = 2226.9583  # BSM I134 (was 2374.6667)
$ws.Cells.Item(134, 10).Value = 0  # BSM J134 (was 2000)
$ws.Cells.Item(134, 11).Value = 6680.874899999999  # BSM K134 (was 7124.000100000001)
$ws.Cells.Item(134, 12).Value = 0  # BSM L134 (was 6000)
$ws.Cells.Item(134, 13).Value = -4145.874899999999  # BSM M134 (was -4589.000100000001)
$ws.Cells.Item(134, 14).Value = $null  # BSM N134 remove (was -11070)

$ws.Cells.Item(141, 8).Value = 75259.664  # BSM H141 (was 77694.5)
$ws.Cells.Item(141, 10).Value = 75259.664  # BSM J141 (was 77694.5)
$ws.Cells.Item(141, 12).Value = 75259.664  # BSM L141 (was 77694.5)
$ws.Cells.Item(141, 14).Value = -85619.664  # BSM N141 (was -88054.5)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 91.125  # CRP H2 (was 317.7143)
$ws.Cells.Item(2, 9).Value = 75.57143000000001  # CRP I2 (was 364.8)
$ws.Cells.Item(2, 11).Value = 75.57143000000001  # CRP K2 (was 364.8)
$ws.Cells.Item(2, 13).Value = 37.42856999999999  # CRP M2 (was -251.8)

$ws.Cells.Item(31, 8).Value = 3037.2222  # CRP H31 (was 2812.2727)
$ws.Cells.Item(31, 10).Value = 2995.6667  # CRP J31 (was 2517.4)
$ws.Cells.Item(31, 12).Value = 2995.6667  # CRP L31 (was 2517.4)
$ws.Cells.Item(31, 14).Value = -3585.6667  # CRP N31 (was -3107.4)

$ws.Cells.Item(34, 8).Value = 3037.2222  # CRP H34 (was 2812.2727)
$ws.Cells.Item(34, 10).Value = 2995.6667  # CRP J34 (was 2517.4)
$ws.Cells.Item(34, 12).Value = 2995.6667  # CRP L34 (was 2517.4)
$ws.Cells.Item(34, 14).Value = -3399.6667  # CRP N34 (was -2921.4)

$ws.Cells.Item(35, 8).Value = 749.75  # CRP H35 (was 1339.8)
$ws.Cells.Item(35, 9).Value = 749.75  # CRP I35 (was 1339.8)
$ws.Cells.Item(35, 11).Value = 749.75  # CRP K35 (was 1339.8)
$ws.Cells.Item(35, 13).Value = -455.75  # CRP M35 (was -1045.8)

$ws.Cells.Item(58, 8).Value = 146438.42  # CRP H58 (was 170176.67)
$ws.Cells.Item(58, 9).Value = 203208.2  # CRP I58 (was 253008)
$ws.Cells.Item(58, 11).Value = 203208.2  # CRP K58 (was 253008)
$ws.Cells.Item(58, 13).Value = -203005.2  # CRP M58 (was -252805)

$ws.Cells.Item(99, 8).Value = 4999  # CRP H99 (was 4332.3335)
$ws.Cells.Item(99, 9).Value = 2998  # CRP I99 (was 2998.5)
$ws.Cells.Item(99, 11).Value = 2998  # CRP K99 (was 2998.5)
$ws.Cells.Item(99, 13).Value = -1500  # CRP M99 (was -1500.5)

$ws.Cells.Item(126, 8).Value = 4999  # CRP H126 (was 4332.3335)
$ws.Cells.Item(126, 9).Value = 2998  # CRP I126 (was 2998.5)
$ws.Cells.Item(126, 11).Value = 8994  # CRP K126 (was 8995.5)
$ws.Cells.Item(126, 13).Value = -6524  # CRP M126 (was -6525.5)

$ws.Cells.Item(134, 8).Value = 47785.363  # CRP H134 (was 43959.418)
$ws.Cells.Item(134, 9).Value = 64262.25  # CRP I134 (was 57330.223)
$ws.Cells.Item(134, 11).Value = 192786.75  # CRP K134 (was 171990.669)
$ws.Cells.Item(134, 13).Value = -190251.75  # CRP M134 (was -169455.669)

$ws.Cells.Item(136, 8).Value = 146438.42  # CRP H136 (was 170176.67)
$ws.Cells.Item(136, 9).Value = 203208.2  # CRP I136 (was 253008)
$ws.Cells.Item(136, 11).Value = 609624.6000000001  # CRP K136 (was 759024)
$ws.Cells.Item(136, 13).Value = -607074.6000000001  # CRP M136 (was -756474)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3000083.2  # CUL H4 (was 2388118.8)
$ws.Cells.Item(4, 9).Value = 2250000  # CUL I4 (was 2728992.8)
$ws.Cells.Item(4, 10).Value = 4500250  # CUL J4 (was 2000)
$ws.Cells.Item(4, 11).Value = 6750000  # CUL K4 (was 8186978.399999999)
$ws.Cells.Item(4, 12).Value = 13500750  # CUL L4 (was 6000)
$ws.Cells.Item(4, 13).Value = -6749888  # CUL M4 (was -8186866.399999999)
$ws.Cells.Item(4, 14).Value = -13500974  # CUL N4 (was -6224)

$ws.Cells.Item(26, 8).Value = 338.5  # CUL H26 (was 303)
$ws.Cells.Item(26, 10).Value = 338.5  # CUL J26 (was 303)
$ws.Cells.Item(26, 12).Value = 1015.5  # CUL L26 (was 909)
$ws.Cells.Item(26, 14).Value = -1591.5  # CUL N26 (was -1485)

$ws.Cells.Item(32, 8).Value = 0  # CUL H32 (was 300)
$ws.Cells.Item(32, 10).Value = 0  # CUL J32 (was 300)
$ws.Cells.Item(32, 12).Value = 0  # CUL L32 (was 900)
$ws.Cells.Item(32, 14).Value = $null  # CUL N32 remove (was -1466)

$ws.Cells.Item(46, 8).Value = 122.5  # CUL H46 (was 99.5)
$ws.Cells.Item(46, 9).Value = 90.77778000000001  # CUL I46 (was 99.5)
$ws.Cells.Item(46, 10).Value = 179.6  # CUL J46 (was 0)
$ws.Cells.Item(46, 11).Value = 272.33334  # CUL K46 (was 298.5)
$ws.Cells.Item(46, 12).Value = 538.8  # CUL L46 (was 0)
$ws.Cells.Item(46, 13).Value = -181.33334  # CUL M46 (was -207.5)
$ws.Cells.Item(46, 14).Value = -720.8  # CUL N46 (was None)

$ws.Cells.Item(55, 8).Value = 5976.625  # CUL H55 (was 6310.2)
$ws.Cells.Item(55, 9).Value = 1461.5  # CUL I55 (was 1624.3334)
$ws.Cells.Item(55, 11).Value = 4384.5  # CUL K55 (was 4873.0002)
$ws.Cells.Item(55, 13).Value = -4207.5  # CUL M55 (was -4696.0002)

$ws.Cells.Item(114, 8).Value = 1727.5  # CUL H114 (was 1426)
$ws.Cells.Item(114, 10).Value = 1727  # CUL J114 (was 1325.3334)
$ws.Cells.Item(114, 12).Value = 5181  # CUL L114 (was 3976.0002)
$ws.Cells.Item(114, 14).Value = -11689  # CUL N114 (was -10484.0002)

$ws.Cells.Item(129, 8).Value = 480614.5  # CUL H129 (was 432838)
$ws.Cells.Item(129, 9).Value = 12391.6  # CUL I129 (was 11437.818)
$ws.Cells.Item(129, 10).Value = 1065893.1  # CUL J129 (was 947882.7)
$ws.Cells.Item(129, 11).Value = 37174.8  # CUL K129 (was 34313.454)
$ws.Cells.Item(129, 12).Value = 3197679.3  # CUL L129 (was 2843648.1)
$ws.Cells.Item(129, 13).Value = -32174.8  # CUL M129 (was -29313.454)
$ws.Cells.Item(129, 14).Value = -3207679.3  # CUL N129 (was -2853648.1)

$ws.Cells.Item(132, 8).Value = 2199.4736  # CUL H132 (was 2458.9285)
$ws.Cells.Item(132, 9).Value = 2698.2856  # CUL I132 (was 3522.25)
$ws.Cells.Item(132, 10).Value = 1908.5  # CUL J132 (was 2033.6)
$ws.Cells.Item(132, 11).Value = 24284.5704  # CUL K132 (was 31700.25)
$ws.Cells.Item(132, 12).Value = 17176.5  # CUL L132 (was 18302.4)
$ws.Cells.Item(132, 13).Value = -21754.5704  # CUL M132 (was -29170.25)
$ws.Cells.Item(132, 14).Value = -22236.5  # CUL N132 (was -23362.4)

$ws.Cells.Item(138, 8).Value = 9999.5  # CUL H138 (was 0)
$ws.Cells.Item(138, 10).Value = 9999.5  # CUL J138 (was 0)
$ws.Cells.Item(138, 12).Value = 29998.5  # CUL L138 (was 0)
$ws.Cells.Item(138, 14).Value = -40278.5  # CUL N138 (was None)

$ws.Cells.Item(139, 8).Value = 0  # CUL H139 (was 503)
$ws.Cells.Item(139, 9).Value = 0  # CUL I139 (was 503)
$ws.Cells.Item(139, 11).Value = 0  # CUL K139 (was 1509)
$ws.Cells.Item(139, 13).Value = $null  # CUL M139 remove (was 3631)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4283.2354  # GSM H102 (was 3805.45)
$ws.Cells.Item(102, 9).Value = 3558.2856  # GSM I102 (was 3194.375)
$ws.Cells.Item(102, 10).Value = 7666.3335  # GSM J102 (was 6249.75)
$ws.Cells.Item(102, 11).Value = 3558.2856  # GSM K102 (was 3194.375)
$ws.Cells.Item(102, 12).Value = 7666.3335  # GSM L102 (was 6249.75)
$ws.Cells.Item(102, 13).Value = -1936.2856  # GSM M102 (was -1572.375)
$ws.Cells.Item(102, 14).Value = -10910.3335  # GSM N102 (was -9493.75)

$ws.Cells.Item(113, 8).Value = 503276.75  # GSM H113 (was 145770.22)
$ws.Cells.Item(113, 9).Value = 504555.5  # GSM I113 (was 102678.9)
$ws.Cells.Item(113, 10).Value = 501998  # GSM J113 (was 253498.5)
$ws.Cells.Item(113, 11).Value = 504555.5  # GSM K113 (was 102678.9)
$ws.Cells.Item(113, 12).Value = 501998  # GSM L113 (was 253498.5)
$ws.Cells.Item(113, 13).Value = -502385.5  # GSM M113 (was -100508.9)
$ws.Cells.Item(113, 14).Value = -506338  # GSM N113 (was -257838.5)

$ws.Cells.Item(122, 8).Value = 4534.2607  # GSM H122 (was 4490.3335)
$ws.Cells.Item(122, 9).Value = 2553  # GSM I122 (was 2592.9285)
$ws.Cells.Item(122, 10).Value = 8249.125  # GSM J122 (was 8285.143)
$ws.Cells.Item(122, 11).Value = 7659  # GSM K122 (was 7778.7855)
$ws.Cells.Item(122, 12).Value = 24747.375  # GSM L122 (was 24855.429)
$ws.Cells.Item(122, 13).Value = -5209  # GSM M122 (was -5328.7855)
$ws.Cells.Item(122, 14).Value = -29647.375  # GSM N122 (was -29755.429)

$ws.Cells.Item(126, 8).Value = 5755.68  # GSM H126 (was 6035.478)
$ws.Cells.Item(126, 9).Value = 4842.8945  # GSM I126 (was 5114.0586)
$ws.Cells.Item(126, 11).Value = 14528.6835  # GSM K126 (was 15342.1758)
$ws.Cells.Item(126, 13).Value = -12058.6835  # GSM M126 (was -12872.1758)

$ws.Cells.Item(132, 8).Value = 255720.5  # GSM H132 (was 503749)
$ws.Cells.Item(132, 9).Value = 338461  # GSM I132 (was 999999)
$ws.Cells.Item(132, 11).Value = 1015383  # GSM K132 (was 2999997)
$ws.Cells.Item(132, 13).Value = -1012853  # GSM M132 (was -2997467)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4081.6667  # LTW H7 (was 4148)
$ws.Cells.Item(7, 9).Value = 4081.6667  # LTW I7 (was 4148)
$ws.Cells.Item(7, 11).Value = 4081.6667  # LTW K7 (was 4148)
$ws.Cells.Item(7, 13).Value = -3969.6667  # LTW M7 (was -4036)

$ws.Cells.Item(22, 8).Value = 37960.805  # LTW H22 (was 48824.207)
$ws.Cells.Item(22, 9).Value = 74826.8  # LTW I22 (was 139674.75)
$ws.Cells.Item(22, 11).Value = 74826.8  # LTW K22 (was 139674.75)
$ws.Cells.Item(22, 13).Value = -74531.8  # LTW M22 (was -139379.75)

$ws.Cells.Item(27, 8).Value = 37960.805  # LTW H27 (was 48824.207)
$ws.Cells.Item(27, 9).Value = 74826.8  # LTW I27 (was 139674.75)
$ws.Cells.Item(27, 11).Value = 74826.8  # LTW K27 (was 139674.75)
$ws.Cells.Item(27, 13).Value = -74719.8  # LTW M27 (was -139567.75)

$ws.Cells.Item(40, 8).Value = 7748.923  # LTW H40 (was 7222.759)
$ws.Cells.Item(40, 9).Value = 5748.6  # LTW I40 (was 5452.727)
$ws.Cells.Item(40, 10).Value = 14416.667  # LTW J40 (was 12785.714)
$ws.Cells.Item(40, 11).Value = 5748.6  # LTW K40 (was 5452.727)
$ws.Cells.Item(40, 12).Value = 14416.667  # LTW L40 (was 12785.714)
$ws.Cells.Item(40, 13).Value = -5612.6  # LTW M40 (was -5316.727)
$ws.Cells.Item(40, 14).Value = -14688.667  # LTW N40 (was -13057.714)

$ws.Cells.Item(82, 8).Value = 2824.5  # LTW H82 (was 2932.8333)
$ws.Cells.Item(82, 9).Value = 2874  # LTW I82 (was 2999)
$ws.Cells.Item(82, 10).Value = 2775  # LTW J82 (was 2866.6667)
$ws.Cells.Item(82, 11).Value = 2874  # LTW K82 (was 2999)
$ws.Cells.Item(82, 12).Value = 2775  # LTW L82 (was 2866.6667)
$ws.Cells.Item(82, 13).Value = -2513  # LTW M82 (was -2638)
$ws.Cells.Item(82, 14).Value = -3497  # LTW N82 (was -3588.6667)

$ws.Cells.Item(85, 8).Value = 2824.5  # LTW H85 (was 2932.8333)
$ws.Cells.Item(85, 9).Value = 2874  # LTW I85 (was 2999)
$ws.Cells.Item(85, 10).Value = 2775  # LTW J85 (was 2866.6667)
$ws.Cells.Item(85, 11).Value = 2874  # LTW K85 (was 2999)
$ws.Cells.Item(85, 12).Value = 2775  # LTW L85 (was 2866.6667)
$ws.Cells.Item(85, 13).Value = -1626  # LTW M85 (was -1751)
$ws.Cells.Item(85, 14).Value = -5271  # LTW N85 (was -5362.6667)

$ws.Cells.Item(122, 8).Value = 4542.136  # LTW H122 (was 4591.8096)
$ws.Cells.Item(122, 10).Value = 5138.2144  # LTW J122 (was 5264.3076)
$ws.Cells.Item(122, 12).Value = 15414.6432  # LTW L122 (was 15792.9228)
$ws.Cells.Item(122, 14).Value = -20314.6432  # LTW N122 (was -20692.9228)

$ws.Cells.Item(126, 8).Value = 4081.6667  # LTW H126 (was 4148)
$ws.Cells.Item(126, 9).Value = 4081.6667  # LTW I126 (was 4148)
$ws.Cells.Item(126, 11).Value = 12245.0001  # LTW K126 (was 12444)
$ws.Cells.Item(126, 13).Value = -9775.000100000001  # LTW M126 (was -9974)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 15499.5  # WVR H4 (was 8749)
$ws.Cells.Item(4, 10).Value = 1999  # WVR J4 (was 1998.6666)
$ws.Cells.Item(4, 12).Value = 1999  # WVR L4 (was 1998.6666)
$ws.Cells.Item(4, 14).Value = -2225  # WVR N4 (was -2224.6666)

$ws.Cells.Item(81, 8).Value = 2251.4443  # WVR H81 (was 2046.6364)
$ws.Cells.Item(81, 9).Value = 2314  # WVR I81 (was 2016.75)
$ws.Cells.Item(81, 11).Value = 4628  # WVR K81 (was 4033.5)
$ws.Cells.Item(81, 13).Value = -3567  # WVR M81 (was -2972.5)

$ws.Cells.Item(84, 8).Value = 2251.4443  # WVR H84 (was 2046.6364)
$ws.Cells.Item(84, 9).Value = 2314  # WVR I84 (was 2016.75)
$ws.Cells.Item(84, 11).Value = 23140  # WVR K84 (was 20167.5)
$ws.Cells.Item(84, 13).Value = -17836  # WVR M84 (was -14863.5)

$ws.Cells.Item(96, 8).Value = 4066  # WVR H96 (was 4010.4443)
$ws.Cells.Item(96, 9).Value = 2586.75  # WVR I96 (was 2461.75)
$ws.Cells.Item(96, 11).Value = 2586.75  # WVR K96 (was 2461.75)
$ws.Cells.Item(96, 13).Value = -1213.75  # WVR M96 (was -1088.75)

$ws.Cells.Item(100, 8).Value = 2377.5  # WVR H100 (was 2202.1428)
$ws.Cells.Item(100, 9).Value = 1853.2  # WVR I100 (was 1736)
$ws.Cells.Item(100, 11).Value = 3706.4  # WVR K100 (was 3472)
$ws.Cells.Item(100, 13).Value = -3165.4  # WVR M100 (was -2931)

$ws.Cells.Item(105, 8).Value = 18699.666  # WVR H105 (was 31999.5)
$ws.Cells.Item(105, 10).Value = 18699.666  # WVR J105 (was 31999.5)
$ws.Cells.Item(105, 12).Value = 18699.666  # WVR L105 (was 31999.5)
$ws.Cells.Item(105, 14).Value = -25687.666  # WVR N105 (was -38987.5)

$ws.Cells.Item(122, 8).Value = 11755  # WVR H122 (was 11793)
$ws.Cells.Item(122, 9).Value = 15505.714  # WVR I122 (was 15560)
$ws.Cells.Item(122, 11).Value = 46517.142  # WVR K122 (was 46680)
$ws.Cells.Item(122, 13).Value = -44067.142  # WVR M122 (was -44230)

$ws.Cells.Item(126, 8).Value = 51772.094  # WVR H126 (was 54501.25)
$ws.Cells.Item(126, 9).Value = 66300.875  # WVR I126 (was 75866.07000000001)
$ws.Cells.Item(126, 10).Value = 5280  # WVR J126 (was 4650)
$ws.Cells.Item(126, 11).Value = 198902.625  # WVR K126 (was 227598.21)
$ws.Cells.Item(126, 12).Value = 15840  # WVR L126 (was 13950)
$ws.Cells.Item(126, 13).Value = -196432.625  # WVR M126 (was -225128.21)
$ws.Cells.Item(126, 14).Value = -20780  # WVR N126 (was -18890)

$ws.Cells.Item(132, 8).Value = 37270.863  # WVR H132 (was 43026.28)
$ws.Cells.Item(132, 9).Value = 41340.652  # WVR I132 (was 48620.863)
$ws.Cells.Item(132, 11).Value = 124021.956  # WVR K132 (was 145862.589)
$ws.Cells.Item(132, 13).Value = -121491.956  # WVR M132 (was -143332.589)
